$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 2792.1333
$ws.Range("I28").Value = 2353
$ws.Range("J28").Value = 3999.75
$ws.Range("K28").Value = 2353
$ws.Range("L28").Value = 3999.75
$ws.Range("M28").Value = -1868
$ws.Range("N28").Value = -4969.75
# Row 98
$ws.Range("H98").Value = 2071.4546
$ws.Range("I98").Value = 2238.6
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 2238.6
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = -740.5999999999999
$ws.Range("N98").Value = -3396
# Row 112
$ws.Range("H112").Value = 1056.4546
$ws.Range("I112").Value = 1080
$ws.Range("J112").Value = 1036.8334
$ws.Range("K112").Value = 3240
$ws.Range("L112").Value = 3110.5002
$ws.Range("M112").Value = -2132
$ws.Range("N112").Value = -5326.5002
# Row 116
$ws.Range("H116").Value = 1860004.4
$ws.Range("I116").Value = 8063.1333
$ws.Range("K116").Value = 8063.1333
$ws.Range("M116").Value = -4621.1333
# Row 122
$ws.Range("H122").Value = 2071.4546
$ws.Range("I122").Value = 2238.6
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 6715.799999999999
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = -4265.799999999999
$ws.Range("N122").Value = -6100
# Row 129
$ws.Range("H129").Value = 2284
$ws.Range("I129").Value = 2969
$ws.Range("K129").Value = 8907
$ws.Range("M129").Value = -3907
# Row 132
$ws.Range("H132").Value = 1078.3396
$ws.Range("I132").Value = 1078.3396
$ws.Range("K132").Value = 3235.0188
$ws.Range("M132").Value = -705.0187999999998

$ws = $wb.Worksheets.Item("ARM")
# Row 15
$ws.Range("H15").Value = 10011
$ws.Range("I15").Value = 10011
$ws.Range("K15").Value = 10011
$ws.Range("M15").Value = -9661
# Row 22
$ws.Range("H22").Value = 3320.3333
$ws.Range("I22").Value = 3880.8
$ws.Range("K22").Value = 3880.8
$ws.Range("M22").Value = -3581.8
# Row 88
$ws.Range("H88").Value = 1149.5
$ws.Range("J88").Value = 1084.5714
$ws.Range("L88").Value = 1084.5714
$ws.Range("N88").Value = -1896.5714
# Row 91
$ws.Range("H91").Value = 1149.5
$ws.Range("J91").Value = 1084.5714
$ws.Range("L91").Value = 1084.5714
$ws.Range("N91").Value = -3892.5714
# Row 122
$ws.Range("H122").Value = 5083.8
$ws.Range("I122").Value = 7459.75
$ws.Range("K122").Value = 22379.25
$ws.Range("M122").Value = -19929.25
# Row 132
$ws.Range("H132").Value = 2200.3125
$ws.Range("I132").Value = 1825.875
$ws.Range("J132").Value = 3323.625
$ws.Range("K132").Value = 5477.625
$ws.Range("L132").Value = 9970.875
$ws.Range("M132").Value = -2947.625
$ws.Range("N132").Value = -15030.875

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 1900
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 2750
$ws.Range("K8").Value = 200
$ws.Range("L8").Value = 2750
$ws.Range("M8").Value = -60
$ws.Range("N8").Value = -3030
# Row 94
$ws.Range("H94").Value = 1966.6666
$ws.Range("J94").Value = 2000
$ws.Range("L94").Value = 2000
$ws.Range("N94").Value = -2902
# Row 105
$ws.Range("H105").Value = 69448.664
$ws.Range("I105").Value = 101849
$ws.Range("K105").Value = 101849
$ws.Range("M105").Value = -100102

$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 93333.336
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
# Row 99
$ws.Range("H99").Value = 4468403.5
$ws.Range("I99").Value = 5108
$ws.Range("J99").Value = 7815875
$ws.Range("K99").Value = 5108
$ws.Range("L99").Value = 7815875
$ws.Range("M99").Value = -3610
$ws.Range("N99").Value = -7818871
# Row 116
$ws.Range("H116").Value = 86830.55499999999
$ws.Range("J116").Value = 86830.55499999999
$ws.Range("L116").Value = 86830.55499999999
$ws.Range("N116").Value = -96008.55499999999
# Row 122
$ws.Range("H122").Value = 2968.3125
$ws.Range("I122").Value = 2024.3
$ws.Range("J122").Value = 4541.6665
$ws.Range("K122").Value = 6072.9
$ws.Range("L122").Value = 13624.9995
$ws.Range("M122").Value = -3622.9
$ws.Range("N122").Value = -18524.9995
# Row 126
$ws.Range("H126").Value = 4468403.5
$ws.Range("I126").Value = 5108
$ws.Range("J126").Value = 7815875
$ws.Range("K126").Value = 15324
$ws.Range("L126").Value = 23447625
$ws.Range("M126").Value = -12854
$ws.Range("N126").Value = -23452565
# Row 132
$ws.Range("H132").Value = 2384.2856
$ws.Range("I132").Value = 2280.0908
$ws.Range("J132").Value = 2766.3333
$ws.Range("K132").Value = 6840.2724
$ws.Range("L132").Value = 8298.999899999999
$ws.Range("M132").Value = -4310.2724
$ws.Range("N132").Value = -13358.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 31.38889
$ws.Range("I38").Value = 14
$ws.Range("J38").Value = 92.25
$ws.Range("K38").Value = 42
$ws.Range("L38").Value = 276.75
$ws.Range("M38").Value = 305
$ws.Range("N38").Value = -970.75
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 107
$ws.Range("H107").Value = 314.4
$ws.Range("I107").Value = 211.5
$ws.Range("J107").Value = 383
$ws.Range("K107").Value = 634.5
$ws.Range("L107").Value = 1149
$ws.Range("M107").Value = 1285.5
$ws.Range("N107").Value = -4989
# Row 131
$ws.Range("H131").Value = 43101.168
$ws.Range("J131").Value = 1891.875
$ws.Range("L131").Value = 5675.625
$ws.Range("N131").Value = -15755.625

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1371.7333
$ws.Range("I102").Value = 1340.6666
$ws.Range("J102").Value = 1496
$ws.Range("K102").Value = 1340.6666
$ws.Range("L102").Value = 1496
$ws.Range("M102").Value = 281.3334
$ws.Range("N102").Value = -4740
# Row 122
$ws.Range("H122").Value = 11212.4
$ws.Range("I122").Value = 16167.167
$ws.Range("K122").Value = 48501.501
$ws.Range("M122").Value = -46051.501
# Row 126
$ws.Range("H126").Value = 4253.4
$ws.Range("I126").Value = 1953
$ws.Range("K126").Value = 5859
$ws.Range("M126").Value = -3389
# Row 132
$ws.Range("H132").Value = 4274.278
$ws.Range("I132").Value = 3204.1333
$ws.Range("K132").Value = 9612.3999
$ws.Range("M132").Value = -7082.3999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 50200.2
$ws.Range("I7").Value = 28500.25
$ws.Range("K7").Value = 28500.25
$ws.Range("M7").Value = -28388.25
# Row 48
$ws.Range("H48").Value = 33333
$ws.Range("I48").Value = 69999
$ws.Range("J48").Value = 15000
$ws.Range("K48").Value = 69999
$ws.Range("L48").Value = 15000
$ws.Range("M48").Value = -69338
$ws.Range("N48").Value = -16322
# Row 55
$ws.Range("H55").Value = 3983.4092
$ws.Range("J55").Value = 9213.625
$ws.Range("L55").Value = 9213.625
$ws.Range("N55").Value = -9559.625
# Row 122
$ws.Range("H122").Value = 25054026
$ws.Range("I122").Value = 84001.39999999999
$ws.Range("K122").Value = 252004.2
$ws.Range("M122").Value = -249554.2
# Row 126
$ws.Range("H126").Value = 50200.2
$ws.Range("I126").Value = 28500.25
$ws.Range("K126").Value = 85500.75
$ws.Range("M126").Value = -83030.75
# Row 132
$ws.Range("H132").Value = 12394.143
$ws.Range("I132").Value = 18842
$ws.Range("K132").Value = 56526
$ws.Range("M132").Value = -53996
# Row 136
$ws.Range("H136").Value = 7143.8096
$ws.Range("I136").Value = 7791.3335
$ws.Range("J136").Value = 6658.1665
$ws.Range("K136").Value = 23374.0005
$ws.Range("L136").Value = 19974.4995
$ws.Range("M136").Value = -20824.0005
$ws.Range("N136").Value = -25074.4995

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 8930944
$ws.Range("I100").Value = 14288914
$ws.Range("J100").Value = 995
$ws.Range("K100").Value = 28577828
$ws.Range("L100").Value = 1990
$ws.Range("M100").Value = -28577287
$ws.Range("N100").Value = -3072
# Row 104
$ws.Range("H104").Value = 1370
$ws.Range("J104").Value = 1370
$ws.Range("L104").Value = 1370
$ws.Range("N104").Value = -8358
# Row 107
$ws.Range("H107").Value = 9576.048000000001
$ws.Range("J107").Value = 9531.421
$ws.Range("L107").Value = 28594.263
$ws.Range("N107").Value = -32434.263
# Row 122
$ws.Range("H122").Value = 2172.7407
$ws.Range("I122").Value = 1895.5264
$ws.Range("J122").Value = 2831.125
$ws.Range("K122").Value = 5686.5792
$ws.Range("L122").Value = 8493.375
$ws.Range("M122").Value = -3236.5792
$ws.Range("N122").Value = -13393.375
# Row 125
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
# Row 126
$ws.Range("H126").Value = 2550
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
# Row 132
$ws.Range("H132").Value = 1611884.4
$ws.Range("I132").Value = 1232.15
$ws.Range("J132").Value = 6213748
$ws.Range("K132").Value = 3696.45
$ws.Range("L132").Value = 18641244
$ws.Range("M132").Value = -1166.45
$ws.Range("N132").Value = -18646304
